# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (column G) values recalculated for rows 2-15
$kValues = @(0, 1, 0, 2, 1, 0, 0, 0, 2, 3, 2, 2, 1, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
